# Daily update at 8 AM UTC
# Row 83 was the previous "latest" row (using the special last-row date
# style). A new row 84 is appended with the next day's data, and the
# special last-row style moves from A83 to A84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the special "latest row" number format currently on A83 before we
# touch it, and the regular date number format used by the rows above it.
$lastRowFormat = $ws.Range("A83").NumberFormat
$normalDateFormat = $ws.Range("A82").NumberFormat

# A83 reverts to the regular date format used by all the other rows above it.
$ws.Range("A83").NumberFormat = $normalDateFormat

# Append the new day's data on row 84.
$ws.Range("A84").Value = 45824
$ws.Range("B84").Value = 356
$ws.Range("C84").Value = 360
$ws.Range("D84").Value = 361

# A84 becomes the new "latest" row and takes on the special last-row number
# format that used to belong to A83.
$ws.Range("A84").NumberFormat = $lastRowFormat
